# Apply "gh-pages output" data refresh to 苏州-漫展信息.xlsx
# Updates "想去人数" (want-to-go count) figures, one refreshed event's
# start-of-range date, and its cover image URL, across the 展览, 演出 and
# 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: 展览 (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F4").Value = 2104
$ws1.Range("F5").Value = 380
$ws1.Range("F6").Value = 663
$ws1.Range("F8").Value = 2094
$ws1.Range("F9").Value = 10869
$ws1.Range("F14").Value = 425
$ws1.Range("E15").Value = "2024.07.19 10:00-07.21 17:00"
$ws1.Range("F15").Value = 9100
$ws1.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202405/Eh06dOvF1715926655440.jpeg"
$ws1.Range("F16").Value = 1119
$ws1.Range("F18").Value = 5327
$ws1.Range("F20").Value = 3375

# ---------------------------------------------------------------------
# Sheet: 演出 (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F2").Value = 25
$ws2.Range("F3").Value = 560

# ---------------------------------------------------------------------
# Sheet: 全部类型 (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value = 2104
$ws4.Range("F5").Value = 380
$ws4.Range("F6").Value = 663
$ws4.Range("F7").Value = 25
$ws4.Range("F9").Value = 2094
$ws4.Range("F10").Value = 560
$ws4.Range("F12").Value = 10869
$ws4.Range("F17").Value = 425
$ws4.Range("E18").Value = "2024.07.19 10:00-07.21 17:00"
$ws4.Range("F18").Value = 9101
$ws4.Range("I18").Value = "//i1.hdslb.com/bfs/openplatform/202405/Eh06dOvF1715926655440.jpeg"
$ws4.Range("F19").Value = 1119
$ws4.Range("F21").Value = 5327
$ws4.Range("F23").Value = 3375
